# Final Prez.pptx - change present tense to past tense in the "pain points"
# bullet list on slide 2 (the shape named "TextBox 7").
#
#   "Uses Excel to " + "store all information" -> one run "Used Excel to store all information"
#   "Has to manually input ..."                 -> "Had to manually input ..."
#   "Cannot track more than 80 items ..."       -> "Could not track more than 80 items ..."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 1 ------------------------------------------------------------
# Today this paragraph is split across two runs ("Uses Excel to " and
# "store all information") and ends with a leftover endParaRPr. Remove the
# whole paragraph - its text AND its trailing paragraph mark - then type a
# brand-new paragraph back in at the same spot. That collapses it into a
# single, normal run (no stray endParaRPr), same as PowerPoint does when a
# paragraph is retyped from scratch.
$run1 = $tr.Find("Uses Excel to ")
$run2 = $tr.Find("store all information")
$wholePara1 = $tr.Characters($run1.Start, $run1.Length + $run2.Length + 1)
[void]$wholePara1.Delete()
[void]$tr.InsertBefore("Used Excel to store all information`r")

# Re-anchor the text range after the structural edit above.
$tr = $shp.TextFrame.TextRange

# --- Paragraph 2 -------------------------------------------------------------
$old2 = $tr.Find("Has to manually input each line of items from a PO or Invoice")
$old2.Text = "Had to manually input each line of items from a PO or Invoice"

# --- Paragraph 3 -------------------------------------------------------------
$old3 = $tr.Find("Cannot track more than 80 items per page, if more than that receiving  must manually create a second log and save it as a separate file")
$old3.Text = "Could not track more than 80 items per page, if more than that receiving  must manually create a second log and save it as a separate file"
